$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 6 existing values (antonym pair -> _0 suffix variants; examples now carry a trailing newline)
$ws.Range("A6").Value = "good_0"
$ws.Range("B6").Value = "bad_0"
$ws.Range("C6").Value = "Eating healthy food is good for you.`n"
$ws.Range("D6").Value = "Smoking is bad for your health."
$ws.Range("E6").Value = "positive"
$ws.Range("F6").Value = "negative`n"

# Add new row 7 - copy formatting from row 6 first, then set values, to reuse existing style index
$ws.Range("A6:F6").Copy()
$ws.Range("A7:F7").PasteSpecial(-4122)

$ws.Range("A7").Value = "good_1"
$ws.Range("B7").Value = "bad_1"
$ws.Range("C7").Value = "Eating healthy food is good for you."
$ws.Range("D7").Value = "Smoking is bad for your health."
$ws.Range("E7").Value = "positive"
$ws.Range("F7").Value = "negative"

# Restore row heights that auto-fit may have disturbed when multi-line text was entered
$ws.Rows.Item(6).RowHeight = 17.25
$ws.Rows.Item(7).RowHeight = 17.25
